# Reorders the database table rows (A2:J16) on Sheet1 into their new
# order, refreshes the dependent formulas (G/H/I) and the "download"
# fill colour that travels with each category block, and finally moves
# the active selection to A4:XFD4 (whole row 4 selected, active cell A4)
# to match the edited file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Target row order (sheet row -> table key). Everything else (D/E/F
#    "From"/"To"/"Frequency") is identical for every row, so it is left
#    untouched.
# ---------------------------------------------------------------------

$ws.Range("A2").Value = "EQ_m_01"
$ws.Range("B2").Value = "Equity Cash Segment - Summary"
$ws.Range("C2").Value = "Secondary Market"
$ws.Range("J2").Value = "Derivatives"

$ws.Range("A3").Value = "EQ_m_02"
$ws.Range("B3").Value = "Equity Cash Segment - Detailed"
$ws.Range("C3").Value = "Secondary Market"
$ws.Range("J3").Value = "Derivatives"

$ws.Range("A4").Value = "ED_m_01"
$ws.Range("B4").Value = "Trends in Equity Derivatives Turnover"
$ws.Range("C4").Value = "Derivatives"
$ws.Range("J4").Value = "Derivatives"

$ws.Range("A5").Value = "MF_m_02"
$ws.Range("B5").Value = "Funds Mobilised by Mutual Funds (Summary)"
$ws.Range("C5").Value = "Mutual Funds"
$ws.Range("J5").Value = "Mutual%20Funds"

$ws.Range("A6").Value = "MF_m_01"
$ws.Range("B6").Value = "Funds Mobilised by Mutual Funds(Detailed MCR)"
$ws.Range("C6").Value = "Mutual Funds"
$ws.Range("J6").Value = "Mutual%20Funds"

$ws.Range("A7").Value = "MF_m_03"
$ws.Range("B7").Value = "Funds Mobilised by Mutual Funds (Equity vs Debt)"
$ws.Range("C7").Value = "Mutual Funds"
$ws.Range("J7").Value = "Mutual%20Funds"

$ws.Range("A8").Value = "MF_m_04"
$ws.Range("B8").Value = "Funds Mobilised by Mutual Funds (Open Vs Close ended)"
$ws.Range("C8").Value = "Mutual Funds"
$ws.Range("J8").Value = "Mutual%20Funds"

$ws.Range("A9").Value = "PM_m_03"
$ws.Range("B9").Value = "Primary Market (Public Equity) IPO, FPO, Rights"
$ws.Range("C9").Value = "Primary Market"
$ws.Range("J9").Value = "Primary%20Market"

$ws.Range("A10").Value = "PM_m_04"
$ws.Range("B10").Value = "Primary Market (Public Equity) - Region Wise"
$ws.Range("C10").Value = "Primary Market"
$ws.Range("J10").Value = "Primary%20Market"

$ws.Range("A11").Value = "PM_m_05"
$ws.Range("B11").Value = "Primary Market (Public Equity) - Sector Wise (Private vs Public)"
$ws.Range("C11").Value = "Primary Market"
$ws.Range("J11").Value = "Primary%20Market"

$ws.Range("A12").Value = "PM_m_06"
$ws.Range("B12").Value = "Primary Market (Public Equity) - Industry Wise"
$ws.Range("C12").Value = "Primary Market"
$ws.Range("J12").Value = "Primary%20Market"

$ws.Range("A13").Value = "PM_m_07"
$ws.Range("B13").Value = "Primary Market (Public Equity) - Size Wise"
$ws.Range("C13").Value = "Primary Market"
$ws.Range("J13").Value = "Primary%20Market"

$ws.Range("A14").Value = "PM_m_08"
$ws.Range("B14").Value = "Primary Market (Public Equity) - Issue Wise (Detailed)"
$ws.Range("C14").Value = "Primary Market"
$ws.Range("J14").Value = "Primary%20Market"

$ws.Range("A15").Value = "PM_m_01"
$ws.Range("B15").Value = "Primary Market - Summary"
$ws.Range("C15").Value = "Primary Market"
$ws.Range("J15").Value = "Primary%20Market"

$ws.Range("A16").Value = "PM_m_02"
$ws.Range("B16").Value = "Primary Market - Detailed"
$ws.Range("C16").Value = "Primary Market"
$ws.Range("J16").Value = "Primary%20Market"

# ---------------------------------------------------------------------
# 2. Refresh the link formulas in G/H/I so they point at their own row.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 16; $r++) {
    $ws.Range("G$r").Formula = "=CONCATENATE(""/database/"",A$r)"
    $ws.Range("H$r").Formula = "=CONCATENATE(""/database/"",A$r,""/download"")"
    $ws.Range("I$r").Formula = "=CONCATENATE(""/database/"",A$r,""/download"")"
}

# ---------------------------------------------------------------------
# 3. The "Download" column fill colour is keyed to each row's category;
#    re-apply it for the new layout (colours unchanged, only which rows
#    they land on).
# ---------------------------------------------------------------------

# Secondary Market / Derivatives rows -> rows 2-4
$ws.Range("H2:H4").Interior.Color = 0xA6AC81
$ws.Range("H2:H4").Interior.PatternColor = 0xFF9999

# Mutual Funds rows -> rows 5-8
$ws.Range("H5:H8").Interior.Color = 0x94E9FF
$ws.Range("H5:H8").Interior.PatternColor = 0xCCFFFF

# Primary Market rows -> rows 9-16
$ws.Range("H9:H16").Interior.Color = 0x6CB6FF
$ws.Range("H9:H16").Interior.PatternColor = 0xCC99FF

# ---------------------------------------------------------------------
# 4. Match the saved selection: whole row 4 selected, active cell A4.
# ---------------------------------------------------------------------
$ws.Range("A4:XFD4").Select()
$excel.ActiveWindow.RangeSelection.Item(1).Activate()
